$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on touched cells so numeric-looking strings
# (e.g. '1.00', '7.20') are preserved verbatim as text, matching the
# original inlineStr cell type instead of being coerced to numbers.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '62.346.23'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -3.15%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.369.60'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -3.97%  '

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '565.48'

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '124.12'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -7.67%  '

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.370.56'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -3.96%  '

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -3.67%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.20'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -5.06%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.119'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -4.79%  '

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -4.80%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.936.10'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -4.21%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.118'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -1.06%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.365.53'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -4.11%  '

$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -6.57%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '62.436.35'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -2.98%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '24.26'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -5.94%  '

$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -7.53%  '

$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -3.25%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.99'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -4.34%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '367.71'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -6.87%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.551'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -4.53%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.501.19'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -4.08%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.13%  '

$ws.Range('B26').NumberFormat = "@"
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').NumberFormat = "@"
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '70.80'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -5.27%  '

$ws.Range('B27').NumberFormat = "@"
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').NumberFormat = "@"
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000104'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -11.33%  '

$ws.Range('B28').NumberFormat = "@"
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').NumberFormat = "@"
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.13%  '

$ws.Range('B29').NumberFormat = "@"
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').NumberFormat = "@"
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.78'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -8.22%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.10'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -7.20%  '

$ws.Range('B31').NumberFormat = "@"
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').NumberFormat = "@"
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.68'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -7.41%  '

$ws.Range('B33').NumberFormat = "@"
$ws.Range('B33').Value = 'RenzoRestakedETH'
$ws.Range('C33').NumberFormat = "@"
$ws.Range('C33').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.396.33'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -4.07%  '

$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -5.67%  '

$ws.Range('B35').NumberFormat = "@"
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').NumberFormat = "@"
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.36'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -6.66%  '

$ws.Range('B36').NumberFormat = "@"
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').NumberFormat = "@"
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '22.48'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -3.78%  '

$ws.Range('B37').NumberFormat = "@"
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').NumberFormat = "@"
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.05'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -6.04%  '

$ws.Range('B38').NumberFormat = "@"
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').NumberFormat = "@"
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '164.31'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -1.83%  '

$ws.Range('B39').NumberFormat = "@"
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').NumberFormat = "@"
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.57'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -5.58%  '

$ws.Range('B40').NumberFormat = "@"
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').NumberFormat = "@"
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.46'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -5.90%  '

$ws.Range('B41').NumberFormat = "@"
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').NumberFormat = "@"
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0743'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -5.84%  '

$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.10%  '

$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.762'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -5.97%  '

$ws.Range('B44').NumberFormat = "@"
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').NumberFormat = "@"
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '41.02'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -2.57%  '

$ws.Range('B45').NumberFormat = "@"
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').NumberFormat = "@"
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '4.19'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -5.71%  '

$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.52'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -8.13%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '22.23'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -11.89%  '

$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.05'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -10.14%  '

$ws.Range('B49').NumberFormat = "@"
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').NumberFormat = "@"
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.57'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -3.65%  '

$ws.Range('B50').NumberFormat = "@"
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').NumberFormat = "@"
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.220.05'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -6.76%  '

$ws.Range('B51').NumberFormat = "@"
$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').NumberFormat = "@"
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.835'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -7.01%  '

